$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "59.784.63"
$ws.Cells.Item(2, 5).Value = "  -0.04%  "

$ws.Cells.Item(3, 4).Value = "2.364.37"
$ws.Cells.Item(3, 5).Value = "  -1.80%  "

$ws.Cells.Item(4, 5).Value = "  +0.02%  "

$ws.Cells.Item(5, 4).Value = "'556.79"
$ws.Cells.Item(5, 5).Value = "  +1.03%  "

$ws.Cells.Item(6, 4).Value = "'133.17"
$ws.Cells.Item(6, 5).Value = "  -2.77%  "

$ws.Cells.Item(7, 5).Value = "  +0.07%  "

$ws.Cells.Item(8, 4).Value = "'0.583"
$ws.Cells.Item(8, 5).Value = "  -1.40%  "

$ws.Cells.Item(9, 5).Value = "  -0.21%  "

$ws.Cells.Item(10, 4).Value = "'5.64"
$ws.Cells.Item(10, 5).Value = "  -0.55%  "

$ws.Cells.Item(11, 5).Value = "  +1.06%  "

$ws.Cells.Item(12, 4).Value = "'0.341"
$ws.Cells.Item(12, 5).Value = "  -3.21%  "

$ws.Cells.Item(13, 4).Value = "'24.28"
$ws.Cells.Item(13, 5).Value = "  -3.83%  "

$ws.Cells.Item(14, 4).Value = "2.790.70"
$ws.Cells.Item(14, 5).Value = "  -1.67%  "

$ws.Cells.Item(15, 4).Value = "59.725.06"
$ws.Cells.Item(15, 5).Value = "  -0.04%  "

$ws.Cells.Item(16, 5).Value = "  +0.02%  "

$ws.Cells.Item(17, 4).Value = "2.379.81"
$ws.Cells.Item(17, 5).Value = "  -1.68%  "

$ws.Cells.Item(18, 4).Value = "'11.06"
$ws.Cells.Item(18, 5).Value = "  -2.12%  "

$ws.Cells.Item(19, 4).Value = "'4.46"
$ws.Cells.Item(19, 5).Value = "  +1.14%  "

$ws.Cells.Item(20, 4).Value = "'319.81"
$ws.Cells.Item(20, 5).Value = "  -2.63%  "

$ws.Cells.Item(21, 4).Value = "'6.64"
$ws.Cells.Item(21, 5).Value = "  -0.38%  "

$ws.Cells.Item(22, 5).Value = "  +0.06%  "

$ws.Cells.Item(23, 4).Value = "'64.21"
$ws.Cells.Item(23, 5).Value = "  -2.83%  "

$ws.Cells.Item(24, 5).Value = "  -0.35%  "

$ws.Cells.Item(25, 5).Value = "  -0.06%  "

$ws.Cells.Item(26, 4).Value = "'8.39"
$ws.Cells.Item(26, 5).Value = "  -2.55%  "

$ws.Cells.Item(27, 4).Value = "'1.37"
$ws.Cells.Item(27, 5).Value = "  -0.05%  "

$ws.Cells.Item(28, 5).Value = "  +1.92%  "

$ws.Cells.Item(29, 4).Value = "0.0₃0755"
$ws.Cells.Item(29, 5).Value = "  -1.73%  "

$ws.Cells.Item(30, 4).Value = "'170.57"
$ws.Cells.Item(30, 5).Value = "  +0.78%  "

$ws.Cells.Item(31, 4).Value = "'6.06"
$ws.Cells.Item(31, 5).Value = "  +0.71%  "

$ws.Cells.Item(32, 4).Value = "'1.11"
$ws.Cells.Item(32, 5).Value = "  +10.33%  "

$ws.Cells.Item(33, 4).Value = "'0.398"
$ws.Cells.Item(33, 5).Value = "  -1.01%  "

$ws.Cells.Item(34, 4).Value = "'18.09"
$ws.Cells.Item(34, 5).Value = "  -2.73%  "

$ws.Cells.Item(35, 2).Value = "USDe"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(35, 4).Value = "'0.999"
$ws.Cells.Item(35, 5).Value = "  +0.01%  "

$ws.Cells.Item(36, 2).Value = "ImmutableX"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).Value = "'1.32"
$ws.Cells.Item(36, 5).Value = "  +0.97%  "

$ws.Cells.Item(37, 5).Value = "  +0.04%  "

$ws.Cells.Item(38, 4).Value = "'4.12"
$ws.Cells.Item(38, 5).Value = "  -1.20%  "

$ws.Cells.Item(39, 4).Value = "'1.58"
$ws.Cells.Item(39, 5).Value = "  -1.14%  "

$ws.Cells.Item(40, 4).Value = "'318.73"
$ws.Cells.Item(40, 5).Value = "  -0.72%  "

$ws.Cells.Item(41, 4).Value = "'38.57"
$ws.Cells.Item(41, 5).Value = "  -2.30%  "

$ws.Cells.Item(42, 4).Value = "'144.76"
$ws.Cells.Item(42, 5).Value = "  +3.52%  "

$ws.Cells.Item(43, 4).Value = "'3.52"
$ws.Cells.Item(43, 5).Value = "  -3.36%  "

$ws.Cells.Item(44, 4).Value = "'0.0965"
$ws.Cells.Item(44, 5).Value = "  -0.10%  "

$ws.Cells.Item(45, 4).Value = "'19.46"
$ws.Cells.Item(45, 5).Value = "  -0.07%  "

$ws.Cells.Item(46, 4).Value = "'0.0509"
$ws.Cells.Item(46, 5).Value = "  -0.83%  "

$ws.Cells.Item(47, 4).Value = "'0.567"
$ws.Cells.Item(47, 5).Value = "  -1.73%  "

$ws.Cells.Item(48, 5).Value = "  -2.63%  "

$ws.Cells.Item(49, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(49, 4).Value = "'11.07"
$ws.Cells.Item(49, 5).Value = "  +0.34%  "

$ws.Cells.Item(50, 5).Value = "  -0.27%  "

$ws.Cells.Item(51, 2).Value = "dogwifhat"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(51, 4).Value = "'1.53"
$ws.Cells.Item(51, 5).Value = "  -1.49%  "
